# Auto-generated edit script: updates cryptos list values (prices/volume %)
# and fixes row ordering for rows 49-51 (TrueUSD/USDD/EnergySwap), per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.485.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.12%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.664.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.60%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4617"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2571"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06132"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.74%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.664.79"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06949"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.54"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.333"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "74.63"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5621"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.85%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.21%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.492.59"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006650"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.95%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.05%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.878.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.55%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.410"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.684"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.187"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.24"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.80"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.370"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "103.70"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.692"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.928"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.52%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07721"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.590"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04263"

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.45%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9396"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5948"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9180"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +11.57%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.50%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "102.09"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01457"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.806"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3683"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.916"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.64%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.05282"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.70%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1099"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.090"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "29.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.26%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.368"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.05%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9982"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.18%  "
